# Apply the latest crypto price/volume snapshot to Sheet1.
# Source data: GitHub Actions scheduled refresh (cryptos.xlsx).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "91.721.49"
$ws.Range("E2").Value = "  +2.34%  "

# Row 3
$ws.Range("D3").Value = "3.163.34"
$ws.Range("E3").Value = "  +2.68%  "

# Row 4
$ws.Range("E4").Value = "  -0.31%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.67"
$ws.Range("E5").Value = "  +1.73%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "621.09"
$ws.Range("E6").Value = "  +0.37%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.13"
$ws.Range("E7").Value = "  +6.91%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.376"
$ws.Range("E8").Value = "  +3.57%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.03%  "

# Row 10
$ws.Range("D10").Value = "3.160.36"
$ws.Range("E10").Value = "  +2.65%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.746"
$ws.Range("E11").Value = "  +4.45%  "

# Row 12
$ws.Range("E12").Value = "  +2.60%  "

# Row 13
$ws.Range("E13").Value = "  -1.37%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.55"
$ws.Range("E14").Value = "  +0.56%  "

# Row 15
$ws.Range("E15").Value = "  +3.76%  "

# Row 16
$ws.Range("D16").Value = "91.255.44"
$ws.Range("E16").Value = "  +2.14%  "

# Row 17
$ws.Range("D17").Value = "3.752.74"
$ws.Range("E17").Value = "  +3.05%  "

# Row 18
$ws.Range("D18").Value = "3.149.25"
$ws.Range("E18").Value = "  +1.85%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.75"
$ws.Range("E19").Value = "  -1.41%  "

# Row 20
$ws.Range("E20").Value = "  +10.63%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.92"
$ws.Range("E21").Value = "  +9.86%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "458.94"
$ws.Range("E22").Value = "  +6.04%  "

# Row 23
$ws.Range("E23").Value = "  -4.66%  "

# Row 24
$ws.Range("E24").Value = "  +5.35%  "

# Row 25
$ws.Range("E25").Value = "  +8.36%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.89"
$ws.Range("E26").Value = "  +8.88%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.06"
$ws.Range("E27").Value = "  +2.79%  "

# Row 28
$ws.Range("D28").Value = "3.318.56"
$ws.Range("E28").Value = "  +2.01%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.29%  "

# Row 30
$ws.Range("E30").Value = "  +40.26%  "

# Row 31
$ws.Range("E31").Value = "  +18.96%  "

# Row 32
$ws.Range("E32").Value = "  +10.27%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.42"
$ws.Range("E33").Value = "  +4.02%  "

# Row 34
$ws.Range("E34").Value = "  +14.21%  "

# Row 35
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.74"
$ws.Range("E35").Value = "  +4.63%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.72"
$ws.Range("E36").Value = "  +8.31%  "

# Row 37
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.932"
$ws.Range("E37").Value = "  -19.04%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "513.31"
$ws.Range("E38").Value = "  +3.51%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.36"
$ws.Range("E39").Value = "  +8.10%  "

# Row 40
$ws.Range("E40").Value = "  +3.24%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.451"
$ws.Range("E41").Value = "  +13.56%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.84"
$ws.Range("E42").Value = "  +6.09%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.46"
$ws.Range("E43").Value = "  -6.13%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.19"
$ws.Range("E44").Value = "  +0.51%  "

# Row 45
$ws.Range("E45").Value = "  -0.09%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "158.77"
$ws.Range("E46").Value = "  +4.70%  "

# Row 47
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.95"
$ws.Range("E47").Value = "  +5.15%  "

# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.712"
$ws.Range("E48").Value = "  +5.48%  "

# Row 49
$ws.Range("E49").Value = "  +6.12%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.50"
$ws.Range("E50").Value = "  +4.18%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.08"
$ws.Range("E51").Value = "  -0.59%  "

Write-Host "Applied all updates"
